$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.858.82"
$ws.Range("E2").Value = "  +4.18%  "
$ws.Range("D3").Value = "3.358.63"
$ws.Range("E3").Value = "  +4.32%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "561.48"
$ws.Range("E5").Value = "  +3.94%  "
$ws.Range("D6").Value = "152.70"
$ws.Range("E6").Value = "  +4.15%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("E10").Value = "  +3.54%  "
$ws.Range("D11").Value = "0.436"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "3.928.56"
$ws.Range("E12").Value = "  +4.11%  "
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").Value = "26.98"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("E15").Value = "  +2.89%  "
$ws.Range("D16").Value = "62.846.60"
$ws.Range("E16").Value = "  +4.09%  "
$ws.Range("D17").Value = "3.355.25"
$ws.Range("E17").Value = "  +5.43%  "
$ws.Range("D18").Value = "6.36"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "13.87"
$ws.Range("E19").Value = "  +4.81%  "
$ws.Range("D20").Value = "386.35"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "8.32"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").Value = "70.27"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  +5.26%  "
$ws.Range("D26").Value = "8.85"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").Value = "0.0₃0954"
$ws.Range("E27").Value = "  +4.54%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "6.59"
$ws.Range("E29").Value = "  +5.41%  "
$ws.Range("E30").Value = "  +3.62%  "
$ws.Range("D31").Value = "5.62"
$ws.Range("E31").Value = "  +3.01%  "
$ws.Range("D32").Value = "22.94"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").Value = "1.30"
$ws.Range("E33").Value = "  +6.40%  "
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Value = "160.29"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +7.75%  "
$ws.Range("D37").Value = "1.88"
$ws.Range("E37").Value = "  +11.28%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").Value = "26.82"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.0744"
$ws.Range("E39").Value = "  +4.73%  "
$ws.Range("D40").Value = "2.828.62"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").Value = "0.0310"
$ws.Range("E41").Value = "  +7.01%  "
$ws.Range("D42").Value = "0.749"
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("D43").Value = "40.65"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "4.26"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("E45").Value = "  +2.97%  "
$ws.Range("D46").Value = "3.395.43"
$ws.Range("E47").Value = "  +5.12%  "
$ws.Range("D48").Value = "0.103"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("D50").Value = "286.24"
$ws.Range("E50").Value = "  +4.30%  "
$ws.Range("D51").Value = "0.798"
$ws.Range("E51").Value = "  -1.22%  "
